$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.693.84"
$ws.Range("E2").Value = "  -2.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.093.83"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.38"
$ws.Range("E5").Value = "  -1.98%  "
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5153"
$ws.Range("E7").Value = "  -1.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4376"
$ws.Range("E8").Value = "  -3.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.92"
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09237"
$ws.Range("E10").Value = "  +1.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.162"
$ws.Range("E11").Value = "  -2.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.75"
$ws.Range("E12").Value = "  -2.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.107.57"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.221"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.746"
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.48"
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001149"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.70"
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06656"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.177"
$ws.Range("E22").Value = "  -2.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.752.16"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.45"
$ws.Range("E24").Value = "  -2.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.318"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.349.01"
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.92"
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.512"
$ws.Range("E28").Value = "  -3.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.03"
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.89"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.137"
$ws.Range("E31").Value = "  -6.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1048"
$ws.Range("E32").Value = "  -2.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.622"
$ws.Range("E33").Value = "  -4.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.176"
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.953"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.302"
$ws.Range("E36").Value = "  +2.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.20"
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02571"
$ws.Range("E38").Value = "  -2.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.7091"
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06702"
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.330"
$ws.Range("E41").Value = "  +4.77%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.40"
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2223"
$ws.Range("E43").Value = "  -5.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6984"
$ws.Range("E44").Value = "  +7.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.28"
$ws.Range("E45").Value = "  -2.72%  "
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.309"
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.621"
$ws.Range("E48").Value = "  -2.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000352"
$ws.Range("E49").Value = "  -5.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.217"
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.71"
$ws.Range("E51").Value = "  -2.27%  "
